$wb = $excel.ActiveWorkbook

# --- Sheet "JatHarcos" (sheet1): two new enemy rows appended after row 31 ---
$wsHarcos = $wb.Worksheets.Item("JatHarcos")

$wsHarcos.Cells.Item(32, 1).Value = "'sfdf"
$wsHarcos.Cells.Item(32, 2).Value = "'Harcos"
$wsHarcos.Cells.Item(32, 3).Value = "'csatabárd"
$wsHarcos.Cells.Item(32, 4).Value = "'/Images/Karakterek/harcos0.png"

$wsHarcos.Cells.Item(33, 1).Value = "'fdfd"
$wsHarcos.Cells.Item(33, 2).Value = "'Harcos"
$wsHarcos.Cells.Item(33, 3).Value = "'lándzsa"
$wsHarcos.Cells.Item(33, 4).Value = "'/Images/Karakterek/harcos0.png"

# --- Sheet "JatMagus" (sheet2): one new enemy row appended after row 28 ---
$wsMagus = $wb.Worksheets.Item("JatMagus")

$wsMagus.Cells.Item(29, 1).Value = "'hghdgssdf"
$wsMagus.Cells.Item(29, 2).Value = "'Mágus"
$wsMagus.Cells.Item(29, 3).Value = "'tűz botja"
$wsMagus.Cells.Item(29, 4).Value = "'/Images/Karakterek/magus1.png"
